# Add test data row to the TTD upload worksheet.
#
# The existing row 3 (Premium Range Purchasers / 31804) and row 4
# (TV Channels Watched Live / 32048) are swapped, and the new top row
# (now row 3) gets a fresh "taxoapitest" segment name plus an updated
# Segment ID (32048 in both A/B would have been a straight copy, but the
# authored data instead stores the new free-text name in column B and
# resets the Price to 0 with default formatting). Row 4 keeps the old
# row-3 values, except column B is normalized to match column A (31804).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTD")

# --- New row 3: TV Channels Watched Live segment, now with a test name ---
$ws.Cells.Item(3, 1).Value = 32048
$ws.Cells.Item(3, 2).Value = "taxoapitest"
$ws.Cells.Item(3, 3).Value = "TV Channels Watched Live (Last 30 Days)"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "Not Buyable"
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = "Media > TV And Film > TV Channels Watched Live (Last 30 Days)"

# --- New row 4: Premium Range Purchasers segment (was row 3) ---
$ws.Cells.Item(4, 1).Value = 31804
$ws.Cells.Item(4, 2).Value = 31804
$ws.Cells.Item(4, 3).Value = "Premium Range Purchasers"
$ws.Cells.Item(4, 4).Value = "Users who prefer premium branded ranges over supermarket own brand"
$ws.Cells.Item(4, 5).Value = "Buyable"
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Cells.Item(4, 7).Value = 1.5
$ws.Cells.Item(4, 8).Value = "UK Kantar Media TGI > Grocery Shopping > Premium Range Purchasers"

# Make TTD the active sheet/tab with the selection left on H9, matching
# the author switching back to this worksheet after adding the test row.
$ws.Activate()
$ws.Range("H9").Select()
